# Trade #2 closed at 2026-02-17 07:57:16 - unknown UNKNOWN +0.000%
#
# Updates the live trading results workbook:
#  - Summary sheet: roll the aggregate stats to reflect the new closed trade
#  - Strategy Status sheet: roll the MarketMaking strategy row
#  - All Trades / MarketMaking sheets: append the new trade as row 3

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Summary
# ---------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1200.01   # Current Capital
$summary.Range("B4").Value = 0.01      # Total P&L $
$summary.Range("B5").Value = 0.1       # Total P&L %
$summary.Range("B6").Value = 2         # Total Trades
$summary.Range("B8").Value = 1         # Losing Trades
$summary.Range("B9").Value = 50        # Win Rate %

# ---------------------------------------------------------------
# Strategy Status - MarketMaking row (row 4)
# ---------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C4").Value = 100.01     # Capital
$status.Range("D4").Value = 2          # Trades
$status.Range("E4").Value = 0.01       # P&L $
$status.Range("F4").Value = 0.01       # P&L %
$status.Range("G4").Value = 50         # Win Rate %

# ---------------------------------------------------------------
# Helper: write the new trade row (row 3) into a trade-log sheet.
# The Date column ("2026-02-17") must stay literal text, not get
# auto-converted into a date serial, so the cell is pre-formatted
# as Text before the value is assigned, then the format is cleared
# again so no stray number format lingers on the cell.
# ---------------------------------------------------------------
function Add-Trade2Row($sheet) {
    $sheet.Cells.Item(3, 1).Value = 2

    $dateCell = $sheet.Cells.Item(3, 2)
    $dateCell.NumberFormat = "@"
    $dateCell.Value = "2026-02-17"
    $dateCell.ClearFormats()

    $sheet.Cells.Item(3, 3).Value = "07:57:10"
    $sheet.Cells.Item(3, 4).Value = "MarketMaking"
    $sheet.Cells.Item(3, 5).Value = "DOWN"
    $sheet.Cells.Item(3, 6).Value = 0.84
    $sheet.Cells.Item(3, 7).Value = 0.83
    $sheet.Cells.Item(3, 8).Value = "CLOSED"
    $sheet.Cells.Item(3, 9).Value = -1.1905
    $sheet.Cells.Item(3, 10).Value = -0.01
    $sheet.Cells.Item(3, 11).Value = 100.01
    $sheet.Cells.Item(3, 12).Value = 0
    $sheet.Cells.Item(3, 13).Value = 0
    $sheet.Cells.Item(3, 14).Value = 0.6
    $sheet.Cells.Item(3, 15).Value = "Normal spread capture: 19600 bps"
    $sheet.Cells.Item(3, 16).Value = "early_exit"
    $sheet.Cells.Item(3, 17).Value = 0.13
}

# ---------------------------------------------------------------
# All Trades
# ---------------------------------------------------------------
$allTrades = $wb.Worksheets.Item("All Trades")
Add-Trade2Row $allTrades

# ---------------------------------------------------------------
# MarketMaking (per-strategy trade log)
# ---------------------------------------------------------------
$marketMaking = $wb.Worksheets.Item("MarketMaking")
Add-Trade2Row $marketMaking
